$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# Row 16 currently has blank cells in B:K and M (left over from the prior
# partial import); backfill them with the literal "nan" placeholder text
# used throughout the rest of the sheet for "no value".
"B16","C16","D16","E16","F16","G16","H16","I16","J16","K16","M16" | ForEach-Object {
    $ws.Range($_).Value = "nan"
}

# Add the new service-history row (event #16 on 2025-10-01): card id, date,
# correction note and servicer, leaving the measurement columns blank
# (same shape as the row that was just added above before the backfill).
$ws.Range("A17").Value = "'16"
"B17","C17","D17","E17","F17","G17","H17","I17","J17","K17","M17" | ForEach-Object {
    $ws.Range($_).Value = "'"
}
$ws.Range("L17").Value = "10\1\2025"
$ws.Range("N17").Value = "تم تغيير الجرائد الخلفيه (1_5_8) ومعايره"
$ws.Range("O17").Value = "الخبير"
